$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# OUTREACH EXPERIENCE: "Invited Speaker, ..." line
#   Old: Invited Speaker, Kettle Moraine Evening with Nature and Science On
#        Tap-Minocqua <tab> <5 spaces> Sep and Jan 2019
#   New: Invited Speaker, Wisconsin Science Festival, Kettle Moraine Evening
#        with Nature, Science On Tap<break>Minocqua <tab> <5 spaces>
#        <8 tabs> <11 spaces> Oct 2020, Sep 2019, Jan 2019
# ---------------------------------------------------------------------------

# 1) Expand the talk-title text (entirely inside one run, so formatting -
#    Selawik Semibold - is preserved). Insert a manual line break before
#    "Minocqua" via the ^l find/replace code.
$d.Content.Find.Execute(
    "Kettle Moraine Evening with Nature and Science On Tap-Minocqua",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Wisconsin Science Festival, Kettle Moraine Evening with Nature, Science On Tap^lMinocqua",
    2) | Out-Null

# 2) Insert extra tabs + spaces ahead of the date text using InsertBefore so
#    the new characters inherit the formatting of the run immediately before
#    them (the non-italic Gill Sans Nova tab run) rather than the italic
#    date run.
$dateRange = $d.Content
$dateRange.Find.Execute("Sep and Jan 2019") | Out-Null
$dateRange.InsertBefore("`t`t`t`t`t`t`t`t           ") | Out-Null

# 3) Update the date text itself (still its own italic run).
$d.Content.Find.Execute(
    "Sep and Jan 2019",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Oct 2020, Sep 2019, Jan 2019",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Re-touch the "ONLINE COMMUNICATION EXPERIENCE" heading so the stale
# lastRenderedPageBreak cache marker on that run is dropped (no visible text
# change).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "ONLINE COMMUNICATION EXPERIENCE",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ONLINE COMMUNICATION EXPERIENCE",
    2) | Out-Null
